$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix company/name text values (comma -> period typos) ---
$ws.Range("E141").Value = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH"
$ws.Range("E224").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E232").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix "Importe" column number formatting (Argentine "1.234,56" -> plain "1234.56") ---
# Values must remain stored as TEXT (they were text before the fix too), so we
# temporarily force a text number format while assigning, then restore the default style.
$importeValues = @{
  2 = "58500.00"
  3 = "73500.00"
  4 = "29600.00"
  5 = "24500.00"
  6 = "84000.00"
  7 = "8000.00"
  8 = "7800.00"
  9 = "142967.50"
  10 = "324000.00"
  11 = "866000.00"
  12 = "1204000.00"
  13 = "11100.00"
  14 = "4274.00"
  15 = "46506.00"
  16 = "3779.00"
  17 = "1650.00"
  18 = "1098.00"
  19 = "6580.00"
  20 = "17000.00"
  21 = "6224.48"
  22 = "65840.00"
  23 = "24357.96"
  24 = "13600.00"
  25 = "500.00"
  26 = "857155.44"
  27 = "99881.95"
  28 = "29400.00"
  29 = "1590723.60"
  30 = "91888.00"
  31 = "32146.60"
  32 = "151229.72"
  33 = "14462.00"
  34 = "48036.99"
  35 = "79840.00"
  36 = "50000.00"
  37 = "35500.00"
  38 = "2332.88"
  39 = "8400.00"
  40 = "49450.00"
  41 = "5200.00"
  42 = "18924.00"
  43 = "6301.00"
  44 = "95472.00"
  45 = "235334.00"
  46 = "110475.00"
  47 = "56910.00"
  48 = "6361699.23"
  49 = "124021.30"
  50 = "161174.40"
  51 = "1740.00"
  52 = "8560.50"
  53 = "1520.00"
  54 = "4320.90"
  55 = "171870.00"
  56 = "11881.76"
  57 = "5500.00"
  58 = "566748.00"
  59 = "33774.98"
  60 = "2340.00"
  61 = "6200.00"
  62 = "17150.00"
  63 = "1455.00"
  64 = "8910.00"
  65 = "153576.85"
  66 = "2697.18"
  67 = "10550.00"
  68 = "234.99"
  69 = "175437.78"
  70 = "19900.00"
  71 = "17205.78"
  72 = "5040.08"
  73 = "419.00"
  74 = "1998.00"
  75 = "8425.97"
  76 = "2365.54"
  77 = "1950.00"
  78 = "3000.80"
  79 = "7567.16"
  80 = "9852.00"
  81 = "125880.95"
  82 = "8140.00"
  83 = "3600.00"
  84 = "9000.00"
  85 = "23240.00"
  86 = "17600.00"
  87 = "52500.00"
  88 = "34500.00"
  89 = "245040.00"
  90 = "79560.00"
  91 = "191900.00"
  92 = "11725.64"
  93 = "6428.62"
  94 = "47880.00"
  95 = "7900.00"
  96 = "34050.00"
  97 = "127331.00"
  98 = "32500.00"
  99 = "1100.00"
  100 = "15700.00"
  101 = "2900.00"
  102 = "53856.00"
  103 = "3000.00"
  104 = "13125.00"
  105 = "29860.00"
  106 = "15730.00"
  107 = "4665.00"
  108 = "18700.00"
  109 = "900.00"
  110 = "5460.00"
  111 = "1008.20"
  112 = "20.00"
  113 = "33.00"
  114 = "25887.38"
  115 = "2040.00"
  116 = "3190.00"
  117 = "69901.64"
  118 = "286.00"
  119 = "1520.00"
  120 = "300042.00"
  121 = "2345.00"
  122 = "11190.00"
  123 = "13357.00"
  124 = "5776.00"
  125 = "2962.50"
  126 = "14870.00"
  127 = "12532.00"
  128 = "63820.72"
  129 = "1828.80"
  130 = "7500.00"
  131 = "1125.20"
  132 = "4432.00"
  133 = "19490.00"
  134 = "1497.72"
  135 = "1910.70"
  136 = "75.00"
  137 = "39600.00"
  138 = "25500.00"
  139 = "76040.00"
  140 = "16490.00"
  141 = "110000.00"
  142 = "92966.00"
  143 = "340717.72"
  144 = "6677.64"
  145 = "197087.00"
  146 = "300965.00"
  147 = "20060.00"
  148 = "11600.00"
  149 = "16500.00"
  150 = "5000.00"
  151 = "6500.00"
  152 = "10000.00"
  153 = "12000.00"
  154 = "6000.00"
  155 = "29300.00"
  156 = "17000.00"
  157 = "89800.00"
  158 = "51200.00"
  159 = "12000.00"
  160 = "103500.00"
  161 = "159500.00"
  162 = "28000.00"
  163 = "139500.00"
  164 = "91000.00"
  165 = "1765400.00"
  166 = "40000.00"
  167 = "42000.00"
  168 = "10000.00"
  169 = "45000.00"
  170 = "6000.00"
  171 = "7500.00"
  172 = "15445.20"
  173 = "3200.00"
  174 = "111000.00"
  175 = "79500.00"
  176 = "10040.21"
  177 = "10300.00"
  178 = "456.19"
  179 = "7291.00"
  180 = "11500.00"
  181 = "677.99"
  182 = "14500.00"
  183 = "244319.00"
  184 = "35000.00"
  185 = "9000.00"
  186 = "16500.00"
  187 = "22000.00"
  188 = "22000.00"
  189 = "8000.00"
  190 = "7000.00"
  191 = "12000.00"
  192 = "20000.00"
  193 = "15000.00"
  194 = "8500.00"
  195 = "12000.00"
  196 = "16000.00"
  197 = "10000.00"
  198 = "9000.00"
  199 = "9000.00"
  200 = "4000.00"
  201 = "10000.00"
  202 = "10000.00"
  203 = "5000.00"
  204 = "24000.00"
  205 = "10000.00"
  206 = "10000.00"
  207 = "5000.00"
  208 = "52166.50"
  209 = "16000.00"
  210 = "10000.00"
  211 = "12000.00"
  212 = "4500.00"
  213 = "5500.00"
  214 = "18000.00"
  215 = "10000.00"
  216 = "10000.00"
  217 = "41300.00"
  218 = "3500.00"
  219 = "35000.00"
  220 = "20000.00"
  221 = "8690.00"
  222 = "6100.00"
  223 = "593.78"
  224 = "18824.00"
  225 = "3823.00"
  226 = "1000.00"
  227 = "5655.36"
  228 = "13922.74"
  229 = "46580.72"
  230 = "23342.00"
  231 = "16918.31"
  232 = "3800.00"
  233 = "10900.00"
  234 = "22339.14"
  235 = "2520.00"
  236 = "3640.00"
  237 = "1234.00"
  238 = "33850.00"
  239 = "95725.00"
  240 = "16505.05"
  241 = "80000.00"
  242 = "40000.00"
  243 = "40000.00"
  244 = "40000.00"
  245 = "80000.00"
  246 = "40000.00"
  247 = "55000.00"
  248 = "40000.00"
  249 = "40000.00"
  250 = "80000.00"
  251 = "80000.00"
  252 = "1290.00"
  253 = "106045.52"
  254 = "183348.00"
  255 = "22945900.95"
  256 = "85000.00"
  257 = "130000.00"
  258 = "8000.00"
  259 = "265000.00"
  260 = "290100.00"
  261 = "326400.00"
  262 = "294000.00"
  263 = "302100.00"
  264 = "308100.00"
  265 = "317700.00"
  266 = "552850.00"
  267 = "290100.00"
  268 = "722000.00"
  269 = "705000.00"
  270 = "385500.00"
  271 = "290100.00"
  272 = "290100.00"
  273 = "580200.00"
  274 = "452200.00"
  275 = "586300.00"
  276 = "834500.00"
  277 = "553600.00"
  278 = "866400.00"
  279 = "580200.00"
  280 = "303950.00"
  281 = "140680.00"
  282 = "1201024.04"
  283 = "405000.00"
  284 = "258720.00"
  285 = "577900.00"
  286 = "4000.00"
  287 = "186000.00"
  288 = "3800.00"
  289 = "9500.00"
  290 = "9900.00"
  291 = "1000.00"
  292 = "8000.00"
  293 = "29000.00"
  294 = "119000.00"
  295 = "211251.15"
  296 = "12000.00"
  297 = "224250.00"
  298 = "25950.00"
  299 = "18500.00"
  300 = "4000.00"
  301 = "76700.00"
  302 = "40820.00"
}

$importeRange = $ws.Range("H2:H302")
$importeRange.NumberFormat = "@"

foreach ($row in $importeValues.Keys) {
  $ws.Cells.Item($row, 8).Value = $importeValues[$row]
}

$importeRange.Style = "Normal"

